# Updates the cryptos list (Price + Volume(1h) columns) with refreshed values.
# Values are stored as literal text (matching the source data's inlineStr
# cells) even when they look numeric (e.g. "98.17"), so we force text entry
# with a leading apostrophe and then reset the cell style back to "Normal"
# to avoid leaving a stray text-number-format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextValue "D2" "25.575.01"
Set-TextValue "E2" "  +2.77%  "
Set-TextValue "D3" "1.669.05"
Set-TextValue "E4" "  +0.16%  "
Set-TextValue "D5" "237.35"
Set-TextValue "E5" "  +1.15%  "
Set-TextValue "D7" "0.4779"
Set-TextValue "E7" "  +1.14%  "
Set-TextValue "E8" "  +2.45%  "
Set-TextValue "D9" "0.06168"
Set-TextValue "E9" "  +1.45%  "
Set-TextValue "D10" "1.669.53"
Set-TextValue "E10" "  +2.29%  "
Set-TextValue "D11" "0.06999"
Set-TextValue "E11" "  +0.95%  "
Set-TextValue "E12" "  +1.32%  "
Set-TextValue "D13" "0.5903"
Set-TextValue "E13" "  -3.21%  "
Set-TextValue "E14" "  +0.84%  "
Set-TextValue "E15" "  +3.63%  "
Set-TextValue "E16" "  -0.01%  "
Set-TextValue "D17" "0.9999"
Set-TextValue "E17" "  +0.18%  "
Set-TextValue "D18" "25.566.32"
Set-TextValue "E18" "  +2.68%  "
Set-TextValue "D19" "0.000006739"
Set-TextValue "E19" "  +2.84%  "
Set-TextValue "E20" "  +3.32%  "
Set-TextValue "D21" "1.884.66"
Set-TextValue "E22" "  +2.50%  "
Set-TextValue "D23" "8.793"
Set-TextValue "E23" "  +3.01%  "
Set-TextValue "D24" "5.263"
Set-TextValue "E24" "  +0.58%  "
Set-TextValue "D25" "136.80"
Set-TextValue "E25" "  +2.16%  "
Set-TextValue "E26" "  +2.01%  "
Set-TextValue "D27" "1.383"
Set-TextValue "E27" "  +1.00%  "
Set-TextValue "D28" "1.719"
Set-TextValue "E28" "  +5.39%  "
Set-TextValue "D29" "104.72"
Set-TextValue "E29" "  +2.05%  "
Set-TextValue "D30" "4.000"
Set-TextValue "E30" "  +6.97%  "
Set-TextValue "D31" "0.07871"
Set-TextValue "E31" "  +2.15%  "
Set-TextValue "D32" "3.629"
Set-TextValue "E32" "  +2.84%  "
Set-TextValue "E33" "  +0.88%  "
Set-TextValue "D34" "2.623"
Set-TextValue "E34" "  +0.94%  "
Set-TextValue "D35" "0.9559"
Set-TextValue "E35" "  +4.28%  "
Set-TextValue "D36" "0.6047"
Set-TextValue "E36" "  +4.93%  "
Set-TextValue "D37" "0.9438"
Set-TextValue "E37" "  +15.50%  "
Set-TextValue "D38" "2.569"
Set-TextValue "E38" "  +0.74%  "
Set-TextValue "D39" "0.9997"
Set-TextValue "E39" "  +0.19%  "
Set-TextValue "D40" "1.854"
Set-TextValue "E40" "  +4.84%  "
Set-TextValue "E41" "  -4.43%  "
Set-TextValue "D42" "98.17"
Set-TextValue "E42" "  +1.10%  "
Set-TextValue "E43" "  +2.37%  "
Set-TextValue "D44" "4.885"
Set-TextValue "E44" "  +3.99%  "
Set-TextValue "D45" "0.1120"
Set-TextValue "E45" "  +3.03%  "
Set-TextValue "D46" "6.214"
Set-TextValue "E46" "  +3.59%  "
Set-TextValue "D47" "0.05267"
Set-TextValue "E47" "  +1.40%  "
Set-TextValue "D48" "29.97"
Set-TextValue "E48" "  +1.93%  "
Set-TextValue "D49" "7.415"
Set-TextValue "E49" "  +4.04%  "
Set-TextValue "E50" "  +0.23%  "
Set-TextValue "D51" "1.207"
Set-TextValue "E51" "  +2.71%  "
